$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("items")
$wsReq = $wb.Worksheets.Item("requirements")

$wsReq.Range("C2").Copy()
$ws.Range("C1").PasteSpecial(-4122)
Write-Host "test-A (req->items) done"
